$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first 3 accelerometer samples (old rows 2-4) were dropped from this
# capture; deleting them shifts rows 5-21 up into 2-18 with a cell-range
# delete (xlShiftUp) so formatting/positioning behaves like Excel's own
# "Delete Cells..." command rather than a blind overwrite.
$ws.Range("A2:C4").Delete(-4162)

# 13 additional accelerometer samples captured on May 9th are appended
# after the existing data (old last row 21 -> new rows 19-31).
$ws.Range("A19").Value = -28.48631326624204
$ws.Range("B19").Value = 16.53431056169849
$ws.Range("C19").Value = -18.38939690429892
$ws.Range("A20").Value = -6.227460877207299
$ws.Range("B20").Value = -5.338711252148522
$ws.Range("C20").Value = -15.25327578487011
$ws.Range("A21").Value = -2.195728558021952
$ws.Range("B21").Value = -9.250842510453845
$ws.Range("C21").Value = -7.278287349931375
$ws.Range("A22").Value = 22.6993431820965
$ws.Range("B22").Value = 9.916446260157828
$ws.Range("C22").Value = 8.042439793580312
$ws.Range("A23").Value = 10.05702973692206
$ws.Range("B23").Value = 4.430093637248776
$ws.Range("C23").Value = 3.092578779130952
$ws.Range("A24").Value = -17.00681210844283
$ws.Range("B24").Value = -24.96814476563599
$ws.Range("C24").Value = -13.94656551924312
$ws.Range("A25").Value = 33.74133814101275
$ws.Range("B25").Value = 5.435616121996308
$ws.Range("C25").Value = -12.86083946612057
$ws.Range("A26").Value = -5.288869998599139
$ws.Range("B26").Value = -18.46977404460015
$ws.Range("C26").Value = -9.24161248559119
$ws.Range("A27").Value = -25.29642678267203
$ws.Range("B27").Value = -20.23905866738133
$ws.Range("C27").Value = 10.65897996633644
$ws.Range("A28").Value = 3.936601254763914
$ws.Range("B28").Value = 18.79409311281742
$ws.Range("C28").Value = 10.11689840227166
$ws.Range("A29").Value = -36.19211913115269
$ws.Range("B29").Value = -74.38502144013573
$ws.Range("C29").Value = 47.36963078799701
$ws.Range("A30").Value = -34.86098349654431
$ws.Range("B30").Value = -1.257336098075001
$ws.Range("C30").Value = -17.16005880880736
$ws.Range("A31").Value = -23.05234499425689
$ws.Range("B31").Value = -7.394026237846306
$ws.Range("C31").Value = -17.08371260982228
